# Morning Sept 5 rest of kelp blade analysis from 8.31
# Fill in "Kelp Before (cm^2)" (Q) and "Kelp After (cm^2)" (R) measurements
# for the Sept 5 trials (rows 2-15) and add comments (W) on measurement
# quality for a subset of those rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Kelp consumption")
$ws.Activate()

# --- Kelp Before / Kelp After measurements -------------------------------
$ws.Range("Q2").Value = 129.304
$ws.Range("R2").Value = 115.282

$ws.Range("Q3").Value = 125.295
$ws.Range("R3").Value = 123.61

$ws.Range("Q4").Value = 97.843000000000004
$ws.Range("R4").Value = 70.075000000000003

$ws.Range("Q5").Value = 123.244
$ws.Range("R5").Value = 106.542

$ws.Range("Q6").Value = 113.479
$ws.Range("R6").Value = 84.423000000000002

$ws.Range("Q7").Value = 119.378
$ws.Range("R7").Value = 105.98

$ws.Range("Q8").Value = 90.466999999999999
$ws.Range("R8").Value = 90.677999999999997

$ws.Range("Q9").Value = 114.098
$ws.Range("R9").Value = 114.75

$ws.Range("Q10").Value = 108.271
$ws.Range("R10").Value = 108.896

$ws.Range("Q11").Value = 104.895
$ws.Range("R11").Value = 102.779

$ws.Range("Q12").Value = 118.879
$ws.Range("R12").Value = 90.667000000000002

$ws.Range("Q13").Value = 128.02099999999999
$ws.Range("R13").Value = 100.726

$ws.Range("Q14").Value = 117.839
$ws.Range("R14").Value = 105.361

$ws.Range("Q15").Value = 88.838999999999999
$ws.Range("R15").Value = 87.475999999999999

# --- Comments on the photo measurements (column W) ------------------------
# New shared strings must be introduced in this order so that they line up
# with indices 97,98,99,100 in the saved workbook:
#   97 -> "before" out of focus
#   98 -> "after" includes shadows
#   99 -> "before" out of focus, "after" includes shadows
#  100 -> "before" included some holes, "after" includes holes
$ws.Range("W6").Value = [char]34 + "before" + [char]34 + " out of focus"
$ws.Range("W3").Value = [char]34 + "after" + [char]34 + " includes shadows"
$ws.Range("W11").Value = [char]34 + "before" + [char]34 + " out of focus, " + [char]34 + "after" + [char]34 + " includes shadows"
$ws.Range("W5").Value = [char]34 + "before" + [char]34 + " included some holes, " + [char]34 + "after" + [char]34 + " includes holes"

$ws.Range("W7").Value = [char]34 + "after" + [char]34 + " included hole"
$ws.Range("W10").Value = [char]34 + "after" + [char]34 + " includes shadows"
$ws.Range("W12").Value = [char]34 + "after" + [char]34 + " includes shadows"
$ws.Range("W13").Value = [char]34 + "before" + [char]34 + " out of focus"
$ws.Range("W14").Value = [char]34 + "before" + [char]34 + " out of focus, " + [char]34 + "after" + [char]34 + " includes shadows"
$ws.Range("W15").Value = [char]34 + "before" + [char]34 + " out of focus"

# --- Best-effort view state changes (scroll position / selection) --------
$win = $excel.ActiveWindow
$win.ScrollColumn = 13
$win.ScrollRow = 1
$ws.Range("W8").Select()
